# Add a new worksheet ("new") after Sheet1, containing a second data table
# (same layout as Sheet1, with new/updated effectiveness values), and make
# the new sheet the active tab with C4 selected.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$newSheet.Name = "new"

# Header row
$newSheet.Range("A1").Value = "effectivness"
$newSheet.Range("B1").Value = "Baseline"
$newSheet.Range("C1").Value = "Modified"
$newSheet.Range("D1").Value = "Interleaved"

# Row 2
$newSheet.Range("A2").Value = "C"
$newSheet.Range("B2").Value = 0.52859999999999996
$newSheet.Range("C2").Value = 0.51
$newSheet.Range("D2").Value = 0.57110000000000005

# Row 3
$newSheet.Range("A3").Value = "B"
$newSheet.Range("B3").Value = 0.43630000000000002
$newSheet.Range("C3").Value = 0.42820000000000003
$newSheet.Range("D3").Value = 0.45250000000000001

# Row 4
$newSheet.Range("A4").Value = 6
$newSheet.Range("B4").Value = 0.49769999999999998
$newSheet.Range("C4").Value = 0.50690000000000002
$newSheet.Range("D4").Value = 0.57050000000000001

# Row 5
$newSheet.Range("A5").Value = 5
$newSheet.Range("B5").Value = 0.47099999999999997
$newSheet.Range("C5").Value = 0.49640000000000001
$newSheet.Range("D5").Value = 0.56069999999999998

# Row 6
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = 0.3972
$newSheet.Range("C6").Value = 0.40260000000000001
$newSheet.Range("D6").Value = 0.43240000000000001

# Row 7
$newSheet.Range("A7").Value = 3
$newSheet.Range("B7").Value = 0.40039999999999998
$newSheet.Range("C7").Value = 0.39860000000000001
$newSheet.Range("D7").Value = 0.44550000000000001

# Row 8
$newSheet.Range("A8").Value = 2
$newSheet.Range("B8").Value = 0.37930000000000003
$newSheet.Range("C8").Value = 0.37240000000000001
$newSheet.Range("D8").Value = 0.42299999999999999

# Row 9
$newSheet.Range("A9").Value = 1
$newSheet.Range("B9").Value = 0.3004
$newSheet.Range("C9").Value = 0.30199999999999999
$newSheet.Range("D9").Value = 0.32150000000000001

# Make the new sheet active and select C4, matching the saved view state
$newSheet.Activate()
$newSheet.Range("C4").Select()
